$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in undertime minutes for Jan 2023 (row 515) and Feb 2023 (row 516) ---
# Set row 516 first so the new shared-string entries land in the same order
# as the authored workbook (UT(0-0-55) before UT(0-0-48)).
$ws.Range("B516").Value = "UT(0-0-55)"
$ws.Range("D516").Value = 0.115
$ws.Range("B515").Value = "UT(0-0-48)"
$ws.Range("D515").Value = 0.1

# --- Insert a new leave-card row for a 1-day absence taken 3/1/2023 ---
$ws.Rows(518).Insert()

# Copy formatting for the new row from the row right below it (same table style band)
$ws.Range("A519:K519").Copy()
$ws.Range("A518:K518").PasteSpecial(-4122)

# The REMARKS cell holds a date here, so pick up the date number format from row 517
$ws.Range("K517").Copy()
$ws.Range("K518").PasteSpecial(-4122)

# Grow the table definition to cover the freshly inserted row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A8:K573"))

# Populate the new row's data
$ws.Range("B518").Value = "A(1-0-0)"
$ws.Range("D518").Value = 1
$ws.Range("K518").Value = 44986
$ws.Range("G518").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Keep the calculated-column formula on the (now last) table row in its canonical form
$ws.Range("G573").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- Update the CONVERTION lookup used to convert undertime minutes to days ---
$ws2 = $wb.Worksheets.Item("CONVERTION")
$ws2.Range("F3").Value = 48

# --- Restore the reported selection from the authored file ---
$ws.Range("D525").Select()
